# Update "想去人数" (interest count) figures in column F across sheets
# 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types) as published
# by the latest scrape run (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 501
$wsExpo.Range("F5").Value = 142
$wsExpo.Range("F6").Value = 141
$wsExpo.Range("F7").Value = 294
$wsExpo.Range("F9").Value = 125
$wsExpo.Range("F10").Value = 703
$wsExpo.Range("F13").Value = 295
$wsExpo.Range("F15").Value = 6244
$wsExpo.Range("F18").Value = 135
$wsExpo.Range("F20").Value = 15028
$wsExpo.Range("F22").Value = 260
$wsExpo.Range("F25").Value = 10941
$wsExpo.Range("F26").Value = 715
$wsExpo.Range("F27").Value = 4272
$wsExpo.Range("F28").Value = 209
$wsExpo.Range("F30").Value = 121

# Sheet "演出"
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 341

# Sheet "全部类型"
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 501
$wsAll.Range("F5").Value = 142
$wsAll.Range("F6").Value = 141
$wsAll.Range("F7").Value = 294
$wsAll.Range("F8").Value = 341
$wsAll.Range("F10").Value = 125
$wsAll.Range("F11").Value = 703
$wsAll.Range("F15").Value = 295
$wsAll.Range("F18").Value = 6244
$wsAll.Range("F21").Value = 135
$wsAll.Range("F23").Value = 15028
$wsAll.Range("F25").Value = 260
$wsAll.Range("F28").Value = 10941
$wsAll.Range("F29").Value = 715
$wsAll.Range("F30").Value = 4272
$wsAll.Range("F31").Value = 209
$wsAll.Range("F33").Value = 121
